# Insert 4 new data rows (one new reporting date, 44762) right before the
# existing row 136, shifting all subsequent rows down by 4 (old row 226
# ends up at new row 230, and the sheet's used range grows to A1:R230).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A136:A139").EntireRow.Insert()

# --- New row 136 ---
$ws.Range("A136").Value = 2
$ws.Range("B136").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C136").Value = "Coquimbo"
$ws.Range("D136").Value = 44762
$ws.Range("E136").Value = 4
$ws.Range("F136").Value = 100112013
$ws.Range("G136").Value = "Alcachofa"
$ws.Range("H136").Value = "Argentina(o)"
$ws.Range("I136").Value = "Extra"
$ws.Range("J136").Value = 500
$ws.Range("K136").Value = 11000
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = 11500
$ws.Range("N136").Value = "$/caja 40 unidades"
$ws.Range("O136").Value = "Provincia de Limarí"
$ws.Range("P136").Value = 288
$ws.Range("Q136").Value = 40
$ws.Range("R136").Value = "Hortaliza"

# --- New row 137 ---
$ws.Range("A137").Value = 2
$ws.Range("B137").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C137").Value = "Coquimbo"
$ws.Range("D137").Value = 44762
$ws.Range("E137").Value = 4
$ws.Range("F137").Value = 100112013
$ws.Range("G137").Value = "Alcachofa"
$ws.Range("H137").Value = "Argentina(o)"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 700
$ws.Range("K137").Value = 10000
$ws.Range("L137").Value = 11000
$ws.Range("M137").Value = 10500
$ws.Range("N137").Value = "$/caja 50 unidades"
$ws.Range("O137").Value = "Provincia de Limarí"
$ws.Range("P137").Value = 210
$ws.Range("Q137").Value = 50
$ws.Range("R137").Value = "Hortaliza"

# --- New row 138 ---
$ws.Range("A138").Value = 2
$ws.Range("B138").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C138").Value = "Coquimbo"
$ws.Range("D138").Value = 44762
$ws.Range("E138").Value = 4
$ws.Range("F138").Value = 100112013
$ws.Range("G138").Value = "Alcachofa"
$ws.Range("H138").Value = "Española"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 1000
$ws.Range("K138").Value = 13000
$ws.Range("L138").Value = 14000
$ws.Range("M138").Value = 13500
$ws.Range("N138").Value = "$/caja 30 unidades"
$ws.Range("O138").Value = "Provincia de Limarí"
$ws.Range("P138").Value = 450
$ws.Range("Q138").Value = 30
$ws.Range("R138").Value = "Hortaliza"

# --- New row 139 ---
$ws.Range("A139").Value = 2
$ws.Range("B139").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C139").Value = "Coquimbo"
$ws.Range("D139").Value = 44762
$ws.Range("E139").Value = 4
$ws.Range("F139").Value = 100112013
$ws.Range("G139").Value = "Alcachofa"
$ws.Range("H139").Value = "Madrigal"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 400
$ws.Range("K139").Value = 12000
$ws.Range("L139").Value = 13000
$ws.Range("M139").Value = 12500
$ws.Range("N139").Value = "$/caja 40 unidades"
$ws.Range("O139").Value = "Provincia del Elquí"
$ws.Range("P139").Value = 312
$ws.Range("Q139").Value = 40
$ws.Range("R139").Value = "Hortaliza"
